$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "Dr. Gehan Adel, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Servinaz Sayed Mohammad, Administrator"
$ws.Range("G3").Value = "Dr. Veronia Rafat, Dr. Eman Tantawi, Dr. Majorelle Magdy, Dr. Asmaa Reda, Administrator, Dr. Hend Mahmoud"
$ws.Range("G4").Value = "Dr. Gehan Adel, Dr. Eman Tantawi, Dr. Majorelle Magdy, Dr. Servinaz Sayed Mohammad, Dr. Asmaa Reda, Dr. Hend Mahmoud"
$ws.Range("G5").Value = "Dr. Amira Sobhy, Dr. Veronia Rafat, Dr. Eman Tantawi, Dr. Asmaa Reda"
$ws.Range("G6").Value = "Dr. Manar Montaser, Dr. Menna tuâ€™Allah Medhat, Dr. Alshimaa Atef, Dr. Majorelle Magdy, Dr. Mohammad El-Tanany"
$ws.Range("G7").Value = "Dr. Menna tu'Alllah Mohammad, Dr. Amera Ahmad Saad, Dr. Lamiaa Ossama, Dr. Nada Mohammad, Dr. Abeer Ragab, Dr. Kerelos Zareef, Dr. Fatma Elhady"
$ws.Range("G8").Value = "Dr. Nada Mohammad, Dr. Abeer Ragab"
$ws.Range("G11").Value = "Dr. Amal Awwad, Dr. Aya Saeed, Dr. Safa Hany"
$ws.Range("G12").Value = "Dr. Eman M. Abo-Sakaya, Dr. Marina Youhanna, Dr. Dina Adel, Dr. Madeha Saeed, Dr. Yasmeena Fattoh, Dr. Amira Ibrahim"
$ws.Range("G13").Value = "Dr. Yasmeena Fattoh, Dr. Esraa Mostafa, Dr. Amira Ibrahim"
$ws.Range("G15").Value = "Dr. Mohammad Safwat, Dr. Rania Ahmad Youssef"
$ws.Range("G17").Value = "Dr. Mohammad Safwat, Dr. Esraa Samy"
$ws.Range("G20").Value = "Dr. Mohammad Safwat, Dr. Mariam Toma Gerges"
$ws.Range("G27").Value = "Dr. Nourham Mostafa, Dr. Hana Amr"
$ws.Range("G30").Value = "Dr. Shorok Mohammad, Dr. Aya Hanafy, Dr. Yassmen Ahmad, Dr. Wafaa Ebida"
